$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to be bumped by
# one day (46061 -> 46062) for every data row (rows 2 through 84).
for ($row = 2; $row -le 84; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
